$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column A values remain text (matches original inlineStr typing),
# so numeric-looking labels like "0","1","2" do not get coerced to numbers.
$ws.Range("A2:A35").NumberFormat = "@"

$ws.Range("A2").Value = "org_day0_fio2"
$ws.Range("B2").Value = 0.1609525404196093
$ws.Range("A3").Value = "org_day0_map"
$ws.Range("B3").Value = 0.08342418131879123
$ws.Range("A4").Value = "2"
$ws.Range("B4").Value = 0.05531493776990871
$ws.Range("A5").Value = "demo_ethnicity___8"
$ws.Range("B5").Value = 0.05402338763611596
$ws.Range("A6").Value = "ord_day0_gcs"
$ws.Range("B6").Value = 0.052952351962667
$ws.Range("A7").Value = "demo_ethnicity___2"
$ws.Range("B7").Value = 0.04232845218768505
$ws.Range("A8").Value = "1"
$ws.Range("B8").Value = 0.04142244898512395
$ws.Range("A9").Value = "0"
$ws.Range("B9").Value = 0.03875800848063816
$ws.Range("A10").Value = "4"
$ws.Range("B10").Value = 0.03269581157771529
$ws.Range("A11").Value = "bl_lab_ast"
$ws.Range("B11").Value = 0.03093333135846039
$ws.Range("A12").Value = "bl_sao2"
$ws.Range("B12").Value = 0.03068561465677581
$ws.Range("A13").Value = "demo_age_years"
$ws.Range("B13").Value = 0.0258637597766083
$ws.Range("A14").Value = "bl_lab_wbc"
$ws.Range("B14").Value = 0.02395425395618913
$ws.Range("A15").Value = "org_day0_platlet"
$ws.Range("B15").Value = 0.02379625561909436
$ws.Range("A16").Value = "5"
$ws.Range("B16").Value = 0.02372321490539694
$ws.Range("A17").Value = "3"
$ws.Range("B17").Value = 0.02286923726919101
$ws.Range("A18").Value = "demo_ethnicity___6"
$ws.Range("B18").Value = 0.02204655834624651
$ws.Range("A19").Value = "bl_kg"
$ws.Range("B19").Value = 0.02058704102466735
$ws.Range("A20").Value = "org_day0_bilirubin"
$ws.Range("B20").Value = 0.01858959588703498
$ws.Range("A21").Value = "bl_lab_inr"
$ws.Range("B21").Value = 0.01698046465775342
$ws.Range("A22").Value = "demo_ethnicity___7"
$ws.Range("B22").Value = 0.01602358306193647
$ws.Range("A23").Value = "bl_temp"
$ws.Range("B23").Value = 0.01602075263325828
$ws.Range("A24").Value = "bl_lab_haemo"
$ws.Range("B24").Value = 0.0157532379245426
$ws.Range("A25").Value = "bl_lab_troponin"
$ws.Range("B25").Value = 0.01532189077858123
$ws.Range("A26").Value = "6"
$ws.Range("B26").Value = 0.01512677680456213
$ws.Range("A27").Value = "bl_hr"
$ws.Range("B27").Value = 0.01503121863517162
$ws.Range("A28").Value = "bl_lab_creatinine"
$ws.Range("B28").Value = 0.0150150594657461
$ws.Range("A29").Value = "bl_lab_alt"
$ws.Range("B29").Value = 0.01473121246954862
$ws.Range("A30").Value = "demo_ethnicity___3"
$ws.Range("B30").Value = 0.01377286872280545
$ws.Range("A31").Value = "bl_resp_rate"
$ws.Range("B31").Value = 0.01270276081363827
$ws.Range("A32").Value = "co_smoking"
$ws.Range("B32").Value = 0.01206909111160555
$ws.Range("A33").Value = "demo_ethnicity___5"
$ws.Range("B33").Value = 0.008742246367801926
$ws.Range("A34").Value = "bl_oxy_status"
$ws.Range("B34").Value = 0.005234932388247297
$ws.Range("A35").Value = "co_dementia___1"
$ws.Range("B35").Value = 0.002552921026881584
